# Append a new row (51) of portfolio data to Sheet1, extending the
# existing table (A1:D50 -> A1:D51) with the 2025-10-05 entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe so Excel stores the date-like text as a literal
# string (matching the existing "Date" column cells) instead of auto
# converting it to a date serial number, then reset the style back to the
# default "Normal" cell style so no stray number-format is left behind.
$ws.Range("A51").Value = "'2025-10-05"
$ws.Range("A51").Style = "Normal"

$ws.Range("B51").Value = 54.45999908447266
$ws.Range("C51").Value = 716.0999755859375
$ws.Range("D51").Value = 328.4500122070312
